# feat: add 2022-Q4 data
#
# Inserts a new "2022-Q4" worksheet right after "总计", carrying the new
# quarter's fund-holding detail (4 rows), and updates the "总计" summary
# sheet to add the new quarter's row plus a new trailing row that restates
# the previously-last "2020-Q4" total (the other quarter sheets just shift
# position - their own content is untouched).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Insert the new "2022-Q4" sheet right after "总计" (i.e. before the
#    sheet that is currently "2021-Q3").
# ---------------------------------------------------------------------
$zj = $wb.Worksheets.Item("总计")
$oldQ3 = $wb.Worksheets.Item("2021-Q3")

$newSheet = $wb.Worksheets.Add($null, $zj)
$newSheet.Name = "2022-Q4"

# Borrow the header / first-column formatting from the sheet that used to
# carry the "2021-Q3" numbers (same look-and-feel across every quarter tab).
$oldQ3.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$oldQ3.Range("A2").Copy()
$newSheet.Range("A2:A5").PasteSpecial(-4122)

# Header row (note: "基金金额" -> "基金规模" on this sheet only).
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Fund rows. Codes / percentages are textual in this workbook (match the
# existing quarter sheets), so a leading apostrophe keeps them as text
# instead of auto-coercing to numbers (dropping leading zeros etc).
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'016370"
$newSheet.Range("C2").Value = "信澳业绩驱动混合A"
$newSheet.Range("D2").Value = "'0.77"
$newSheet.Range("E2").Value = "'30.31"
$newSheet.Range("F2").Value = "'4.43"
$newSheet.Range("G2").Value = "'0.0341"
$newSheet.Range("H2").Value = 1

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "'011155"
$newSheet.Range("C3").Value = "金鹰责任投资混合A"
$newSheet.Range("D3").Value = "'0.71"
$newSheet.Range("E3").Value = "'92.75"
$newSheet.Range("F3").Value = "'2.85"
$newSheet.Range("G3").Value = "'0.0202"
$newSheet.Range("H3").Value = 9

$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "'016371"
$newSheet.Range("C4").Value = "信澳业绩驱动混合C"
$newSheet.Range("D4").Value = "'0.25"
$newSheet.Range("E4").Value = "'30.31"
$newSheet.Range("F4").Value = "'4.43"
$newSheet.Range("G4").Value = "'0.0111"
$newSheet.Range("H4").Value = 1

$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = "'011156"
$newSheet.Range("C5").Value = "金鹰责任投资混合C"
$newSheet.Range("D5").Value = "'0.21"
$newSheet.Range("E5").Value = "'92.75"
$newSheet.Range("F5").Value = "'2.85"
$newSheet.Range("G5").Value = "'0.0060"
$newSheet.Range("H5").Value = 9

# ---------------------------------------------------------------------
# 2) Update the "总计" summary sheet: row 2 becomes the new 2022-Q4
#    total, the old rows 2-5 shift down to rows 3-6, and a brand new
#    trailing row is appended for 2020-Q4.
# ---------------------------------------------------------------------
$zj.Range("A6").Value = 4
$zj.Range("B6").Value = "2020-Q4"
$zj.Range("C6").Value = 2
$zj.Range("D6").Value = 0.12
$zj.Range("A6").NumberFormat = $zj.Range("A5").NumberFormat
$zj.Range("A5").Copy()
$zj.Range("A6").PasteSpecial(-4122)
$zj.Range("A6").Value = 4

$zj.Range("B5").Value = "2021-Q1"
$zj.Range("D5").Value = 0.23

$zj.Range("B4").Value = "2021-Q2"
$zj.Range("D4").Value = 0.22

$zj.Range("B3").Value = "2021-Q3"
$zj.Range("D3").Value = 0.13

$zj.Range("B2").Value = "2022-Q4"
$zj.Range("C2").Value = 4
$zj.Range("D2").Value = 0.07
